$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column at P (old column P "AnnualizedMoM-CPI-Inflation" shifts to Q)
$ws.Columns("P:P").Insert()

# Insert a new row at 16 (old row 16 "AnnualizedMoM-CPI-Inflation" shifts to 17)
$ws.Rows("16:16").Insert()

# New header for column P / row 16
$ws.Range("P1").Value2 = "RentalPriceAvg%Change"
$ws.Range("A16").Value2 = "RentalPriceAvg%Change"

# Copy the header style (bold + border) from neighboring header cells so the
# newly inserted header cells match the look of the rest of the table.
$ws.Range("O1").Copy() | Out-Null
$ws.Range("P1").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$ws.Range("A15").Copy() | Out-Null
$ws.Range("A16").PasteSpecial(-4122) | Out-Null  # xlPasteFormats
$excel.CutCopyMode = $false

# New correlation values for column P (rows 2-15), RentalPriceAvg%Change vs each other variable
$ws.Range("P2").Value2 = 0.05889156954060009
$ws.Range("P3").Value2 = -0.04067178712249819
$ws.Range("P4").Value2 = -0.05064878077826373
$ws.Range("P5").Value2 = 0.02935561778095505
$ws.Range("P6").Value2 = 0.1415789404803691
$ws.Range("P7").Value2 = -0.02105499089153281
$ws.Range("P8").Value2 = 0.1579867334488559
$ws.Range("P9").Value2 = -0.1071148522496317
$ws.Range("P10").Value2 = -0.18170706190952
$ws.Range("P11").Value2 = -0.07659417116168638
$ws.Range("P12").Value2 = 0.03118207290890015
$ws.Range("P13").Value2 = -0.03672035828256448
$ws.Range("P14").Value2 = 0.04763660614841278
$ws.Range("P15").Value2 = 0.3399132241931659

# New correlation values for row 16 (cols B-O), mirrors column P (symmetric matrix)
$ws.Range("B16").Value2 = 0.05889156954060009
$ws.Range("C16").Value2 = -0.04067178712249819
$ws.Range("D16").Value2 = -0.05064878077826373
$ws.Range("E16").Value2 = 0.02935561778095505
$ws.Range("F16").Value2 = 0.1415789404803691
$ws.Range("G16").Value2 = -0.02105499089153281
$ws.Range("H16").Value2 = 0.1579867334488559
$ws.Range("I16").Value2 = -0.1071148522496317
$ws.Range("J16").Value2 = -0.18170706190952
$ws.Range("K16").Value2 = -0.07659417116168638
$ws.Range("L16").Value2 = 0.03118207290890015
$ws.Range("M16").Value2 = -0.03672035828256448
$ws.Range("N16").Value2 = 0.04763660614841278
$ws.Range("O16").Value2 = 0.3399132241931659

# Diagonal: RentalPriceAvg%Change correlated with itself
$ws.Range("P16").Value2 = 1

# Intersection of RentalPriceAvg%Change (row16/col P) with AnnualizedMoM-CPI-Inflation (row17/col Q)
$ws.Range("Q16").Value2 = -0.02500055734430478
$ws.Range("P17").Value2 = -0.02500055734430478

# Ensure the workbook's used dimension is recalculated properly
$ws.Calculate()
